$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.661.45"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.308.36"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'579.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.58%  "
$ws.Range("D6").Value = "'184.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.303.73"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").Value = "'0.178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").Value = "'46.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "'635.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "3.835.63"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "65.841.04"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "'17.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D20").Value = "3.303.47"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "'11.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "'0.891"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "'17.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").Value = "'100.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'5.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'9.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("D29").Value = "'30.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "'8.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'6.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'594.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.36%  "
$ws.Range("D33").Value = "'3.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.27%  "
$ws.Range("D34").Value = "'10.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").Value = "3.848.16"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'55.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "0.0₃0702"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("D42").Value = "'32.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.44%  "
$ws.Range("D43").Value = "'3.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.21%  "
$ws.Range("D44").Value = "'2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.70%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'0.0409"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("D47").Value = "'3.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.39%  "
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'2.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "'130.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.28%  "
